{"js": "// Office.js (Word JavaScript API) script implementing:\n//  1. Insert a new \"MovingIn, \" run at the start of the KitStatus enum\n//     value list paragraph (the one beginning \"AwaitingPickup, ...\").\n//  2. Relocate the \"_GoBack\" bookmark so it sits immediately after the\n//     newly-inserted \"MovingIn, \" text (it used to live alone in the\n//     empty paragraph right after \"List<Kit> KitsOnStand\").\n//  3. Remove the stray \"DoXXX \" token from the\n//     \"//This agent has no associated //DoXXX animations\" comment so it\n//     reads \"//This agent has no associated //animations\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their (unique) text content.\nlet enumValuesParaIndex = -1;\nlet animationsParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (enumValuesParaIndex === -1 && t.indexOf(\"AwaitingPickup\") === 0) {\n    enumValuesParaIndex = i;\n  }\n  if (animationsParaIndex === -1 && t.indexOf(\"This agent has no associated\") !== -1) {\n    animationsParaIndex = i;\n  }\n}\n\nif (enumValuesParaIndex === -1) {\n  throw new Error(\"Could not find the KitStatus enum value paragraph.\");\n}\nif (animationsParaIndex === -1) {\n  throw new Error(\"Could not find the '//This agent...' paragraph.\");\n}\n\nconst enumValuesPara = paragraphs.items[enumValuesParaIndex];\nconst animationsPara = paragraphs.items[animationsParaIndex];\n\n// --- 1. Insert \"MovingIn, \" at the very start of the enum paragraph ---\nconst insertStart = enumValuesPara.getRange(\"Start\");\ninsertStart.insertText(\"MovingIn, \", \"Start\");\nawait context.sync();\n\n// --- 2. Move the \"_GoBack\" bookmark to right after \"MovingIn, \" ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst movingInMatches = enumValuesPara.search(\"MovingIn, \", { matchCase: true });\nmovingInMatches.load(\"items\");\nawait context.sync();\n\nconst bookmarkSpot = movingInMatches.items[0].getRange(\"End\");\nbookmarkSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 3. Strip \"DoXXX \" out of the animations comment paragraph ---\nconst doxxxMatches = animationsPara.search(\"DoXXX \", { matchCase: true });\ndoxxxMatches.load(\"items\");\nawait context.sync();\n\ndoxxxMatches.items[0].insertText(\"\", \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script implementing:\n#  1. Insert a new \"MovingIn, \" run at the start of the KitStatus enum\n#     value list paragraph (the one beginning \"AwaitingPickup, ...\").\n#  2. Relocate the \"_GoBack\" bookmark so it sits immediately after the\n#     newly-inserted \"MovingIn, \" text (it used to live alone in the\n#     empty paragraph right after \"List<Kit> KitsOnStand\").\n#  3. Remove the stray \"DoXXX \" token from the\n#     \"//This agent has no associated //DoXXX animations\" comment so it\n#     reads \"//This agent has no associated //animations\".\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert \"MovingIn, \" right before \"AwaitingPickup\" ---\n$insertRange = $d.Content\n$insertRange.Find.Text = \"AwaitingPickup\"\n$insertRange.Find.Execute() | Out-Null\n$insertRange.Collapse(1)            # wdCollapseStart\n$insertRange.InsertBefore(\"MovingIn, \")\n\n# --- 2. Move the \"_GoBack\" bookmark to sit right after \"MovingIn, \" ---\n$d.Bookmarks(\"_GoBack\").Delete()\n\n$bookmarkRange = $d.Content\n$bookmarkRange.Find.Text = \"MovingIn, \"\n$bookmarkRange.Find.Execute() | Out-Null\n$bookmarkRange.Collapse(0)          # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# --- 3. Strip \"DoXXX \" out of the animations comment paragraph ---\n$replaceRange = $d.Content\n$replaceRange.Find.Text = \"DoXXX \"\n$replaceRange.Find.Replacement.Text = \"\"\n$replaceRange.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n"}
